# "some basic data for analysis"
#
# The workbook has three sheets: Natural, Grammatical, Failed Catch.
# Natural (sheet1) and Grammatical (sheet2) each had a trailing
# AVERAGE() summary row (row 25 and row 29 respectively) that gets
# removed here. Clearing the row's contents (rather than a full
# row-delete/shift) matches the target: later rows keep their original
# row numbers (e.g. Natural's row 44 stays row 44) while Grammatical's
# used range shrinks from row 29 down to row 26.
#
# The active sheet/selection also flips from Natural (was selected,
# cell K1) to Grammatical (now selected, with the just-cleared summary
# row highlighted as the current selection on each sheet).

$wb = $excel.ActiveWorkbook

$natural = $wb.Worksheets.Item("Natural")
$grammatical = $wb.Worksheets.Item("Grammatical")

# Remove the AVERAGE() summary rows without shifting subsequent rows.
$natural.Rows.Item(25).ClearContents()
$grammatical.Rows.Item(29).ClearContents()

# Leave a selection on the (now empty) former summary row for Natural.
$natural.Activate()
$natural.Range("A25:XFD25").Select()

# Grammatical ends up the active/selected sheet, with its own former
# summary row selected.
$grammatical.Activate()
$grammatical.Range("A29:XFD29").Select()
